$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / request info block ---
# H7 label "Deptartment:" / I7 value
$ws.Range("I7").Value = "IT Department"

# H8 label "Requestor:" / I8 value
$ws.Range("I8").Value = "Stephine David Severino"

# H9 label "Urgency No.:" / I9 value
$ws.Range("I9").Value = 1

# A8 "Date Prepared" value
$ws.Range("C8").Value = [DateTime]"2019-06-26"

# A10 "PR No." value
$ws.Range("C10").Value = "PR-2948-2984"

# Purpose
$ws.Range("C11").Value = "Stator"

# End-Use
$ws.Range("C12").Value = "Warehouse"

# --- Line items table ---
# Row 14
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "kg/s"
$ws.Range("D14").Value = 12111
$ws.Range("E14").Value = "Lumber"
$ws.Range("J14").Value = [DateTime]"2019-02-06"

# Row 15
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "kg/s"
$ws.Range("D15").Value = 1222
$ws.Range("E15").Value = "Nails"
$ws.Range("J15").Value = [DateTime]"2019-02-16"

# Row 16
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "kg/s"
$ws.Range("D16").Value = 1233
$ws.Range("E16").Value = "Hasp"
$ws.Range("J16").Value = [DateTime]"2019-02-18"

# --- Selection / active cell matches the saved view state ---
$ws.Range("C10:E10").Select() | Out-Null
